# Append the 04-11-2025 Gold data row (row 51) to Sheet1, mirroring the
# existing table: col A holds the date label, col B holds the "price of
# gold" paragraph (reused verbatim from the last day that had one, row 44 -
# the same value already used for several other gap days further up).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- A51: write "04-11-2025" as a literal text label -----------------------
# Assigning a date-shaped literal straight to .Value lets Excel's normal
# "smart" entry parsing turn it into a date serial, which is not what the
# source data does (every date in this sheet is plain text). Instead,
# compute it as a text formula result in a scratch cell, then copy only the
# *value* into A51 - this keeps it a plain string cell.
$ws.Range("A100").Formula = '="04-11-2025"'
$ws.Range("A100").Copy()
$ws.Range("A51").PasteSpecial(-4163)   # xlPasteValues
$ws.Range("A100").Clear()
$ws.Rows.Item(100).Delete()            # restore the sheet's used range

# --- B51: reuse the existing price paragraph from row 44 -------------------
$ws.Range("B51").Value = $ws.Range("B44").Text

# --- Match formatting of the preceding data row (row 50) -------------------
$ws.Range("A50:B50").Copy()
$ws.Range("A51:B51").PasteSpecial(-4122)   # xlPasteFormats
